# Rename 'Codelists' to 'Cells' and make it the active sheet/tab
# (Close #256)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Codelists")
$ws.Name = "Cells"

# The active/selected tab moves from "Variables" to the renamed "Cells" sheet
$ws.Activate()
